# Atualização de bases das ligas, do dia: 28-06-2024 às 19:47
#
# This update re-synchronises a handful of match rows whose results/odds
# had been written to the wrong row in the previous refresh. The fix is a
# pure re-shuffle of the data between rows (the "id" running index in
# column A is untouched) within four independent groups of rows:
#   34 <-> 35
#   37 <- 39 <- 38 <- 37   (3-way cycle)
#   40 <-> 42
#   395 <- 396 <- 397 <- 398 <- 395   (4-way cycle)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

function Get-RowValues {
    param($row)
    $vals = @{}
    foreach ($c in $dataCols) {
        $cell = $ws.Range("$c$row")
        if ($cell.Text -eq "") {
            # cell is genuinely absent in the source row - keep that absence
            $vals[$c] = $null
        } else {
            $vals[$c] = $cell.Value2
        }
    }
    return $vals
}

function Set-RowValues {
    param($row, $vals)
    foreach ($c in $dataCols) {
        $cell = $ws.Range("$c$row")
        $v = $vals[$c]
        if ($null -eq $v) {
            $cell.ClearContents() | Out-Null
        } else {
            $cell.Value2 = $v
        }
    }
}

# Snapshot every row in every group BEFORE any of them get overwritten,
# so the permutation is applied against the original values.
$row34 = Get-RowValues 34
$row35 = Get-RowValues 35

$row37 = Get-RowValues 37
$row38 = Get-RowValues 38
$row39 = Get-RowValues 39

$row40 = Get-RowValues 40
$row42 = Get-RowValues 42

$row395 = Get-RowValues 395
$row396 = Get-RowValues 396
$row397 = Get-RowValues 397
$row398 = Get-RowValues 398

# Group 1: swap rows 34 and 35
Set-RowValues 34 $row35
Set-RowValues 35 $row34

# Group 2: 3-way cycle 37 -> 38 -> 39 -> 37
Set-RowValues 37 $row39
Set-RowValues 38 $row37
Set-RowValues 39 $row38

# Group 3: swap rows 40 and 42
Set-RowValues 40 $row42
Set-RowValues 42 $row40

# Group 4: 4-way cycle 395 -> 396 -> 397 -> 398 -> 395
Set-RowValues 395 $row396
Set-RowValues 396 $row397
Set-RowValues 397 $row398
Set-RowValues 398 $row395
